# correção nos dados e inicio da analise PNAD 2009
#
# The sheet originally had several "group header" rows (sexo, cor ou raça,
# grupos de idade, nível de instrução, classes de rendimento mensal
# domiciliar per capita) which carried no data of their own, plus two
# trailing footnote-only rows at the bottom. This edit removes those
# label-only rows (Excel shifts the remaining rows up to close the gaps),
# and fixes the B2 header cell, which incorrectly held the stray label
# "unnamed: 1_level_1" instead of repeating "total".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the mis-labelled header cell (B2): "unnamed: 1_level_1" -> "total"
$ws.Range("B2").Value = "total"

# Remove the now-redundant group-header rows and the trailing footnote
# rows. Deleting from the bottom up keeps the earlier row numbers valid.
$rowsToDelete = @(35, 34, 27, 19, 13, 8, 5)
foreach ($r in $rowsToDelete) {
    $ws.Rows($r).Delete()
}
